$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Coin/Link text for rows with a name change (B and C columns)
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

# Price column (D) values must stay as literal text, not be auto-parsed as numbers/dates.
# Force text format, write the value, then restore the default "Normal" style so no
# stray number-format style attribute is left behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.059.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.834.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6347'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07583'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2949'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07752'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.835.95'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.003'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6705'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009786'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.111'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.091.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.216'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1405'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.536'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.124'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.057'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.207'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05387'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.863'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7509'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.141'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.664'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.237.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.764'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01793'
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9047'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.00000000127'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.986.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5115'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4088'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.074'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05781'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.770'
$ws.Range("D51").Style = "Normal"

# Volume(1h) column (E) values (percentage text with padding spaces)
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("E6").Value = '  +2.04%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +2.92%  '
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("E16").Value = '  +7.39%  '
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  +3.53%  '
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("E33").Value = '  +1.25%  '
$ws.Range("E34").Value = '  +2.26%  '
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("E37").Value = '  -3.63%  '
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("E40").Value = '  +4.76%  '
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("E44").Value = '  +6.04%  '
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("E46").Value = '  +1.79%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  +3.12%  '
$ws.Range("E49").Value = '  +2.64%  '
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("E51").Value = '  +1.72%  '
